$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("bcsstm20.mtx", "MP", 47756316.74668948, 146, 0.01424145698547363, 485),
    @("bcsstm20.mtx", "MP_Aitken", 47811793.47930697, 113, 0.007976293563842773, 485),
    @("bcsstm19.mtx", "MP", 39536719.78976859, 130, 0.02557182312011719, 817),
    @("bcsstm19.mtx", "MP_Aitken", 39549333.10608483, 87, 0.0195920467376709, 817),
    @("bcsstm11.mtx", "MP", 31.57768332328855, 18, 0.02042293548583984, 1473),
    @("bcsstm11.mtx", "MP_Aitken", 31.57822523533776, 11, 0.01172065734863281, 1473),
    @("bcsstm22.mtx", "MP", 0.00973740970842194, 29, 0.001723289489746094, 138),
    @("bcsstm22.mtx", "MP_Aitken", 0.009737699457095711, 17, 0.0009152889251708984, 138),
    @("bcsstm12.mtx", "MP", 13.40996996333015, 122, 0.1330757141113281, 1473),
    @("bcsstm12.mtx", "MP_Aitken", 13.41578254765199, 75, 0.08216667175292969, 1473),
    @("bcsstm07.mtx", "MP", 2509.898102518885, 103, 0.008035898208618164, 420),
    @("bcsstm07.mtx", "MP_Aitken", 2510.150448699714, 66, 0.00531005859375, 420),
    @("bcsstm13.mtx", "MP", 177.2629905522275, 43, 0.08952641487121582, 2003),
    @("bcsstm13.mtx", "MP_Aitken", 177.2664225281739, 22, 0.0450890064239502, 2003),
    @("bcsstm08.mtx", "MP", 1444020.255745148, 29, 0.0133514404296875, 1074),
    @("bcsstm08.mtx", "MP_Aitken", 1444053.928888656, 16, 0.005799770355224609, 1074),
    @("bcsstm09.mtx", "MP", 0.000259020384904, 3, 0.001756668090820312, 1083),
    @("bcsstm09.mtx", "MP_Aitken", 0.000259020384904, 3, 0.002063512802124023, 1083),
    @("bcsstm10.mtx", "MP", 28131.9510377264, 141, 0.06348395347595215, 1086),
    @("bcsstm10.mtx", "MP_Aitken", 28140.67049436026, 78, 0.02787899971008301, 1086)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $ws.Cells.Item($row, 5).Value = $item[3]
    $ws.Cells.Item($row, 6).Value = $item[4]
    $ws.Cells.Item($row, 7).Value = $item[5]
    $row++
}
